# Apply "chapters 15,16 and 17 completed" update to the documentation workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 (date 43977 / 2020-06-15-ish) - chapter 15 completed
$ws.Range("B23").Value = 4
$ws.Range("C23").Value = "chapter 15 completed"
$ws.Range("D23").Value = "had to cheat on the third exercise as after 2.5 hours of work on the prev onesmy head simply did not work at all"

# Row 24 - chapter 16 studied, first two exercises
$ws.Range("B24").Value = 5
$ws.Range("C24").Value = "studied chapter 16, completed first two exercises"
$ws.Range("D24").Value = "had to look up some solutions as mine ones did not perform properly. Also had to look up some explainations on YT as the book did not give enough info"

# Row 25 - finished chapter 16
$ws.Range("B25").Value = 2
$ws.Range("C25").Value = "finished chapter 16"

# Row 26 - free day
$ws.Range("B26").Value = 0
$ws.Range("D26").Value = "freeday"

# Row 27 - chapter 17 completed
$ws.Range("B27").Value = 4
$ws.Range("C27").Value = "chapter 17 completed"

# Update column D width to fit the new longer text (closest reachable value
# to the target raw width of 116.5546875 given the host's sixth-of-a-
# character rounding grid for ColumnWidth).
$ws.Columns("D").ColumnWidth = 115.6

# Update the current view/selection to match the final saved state
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("C27").Select()
